$d = $word.ActiveDocument

$pairs = @(
    @("54×38=", "92×87="),
    @("72×44=", "16×30="),
    @("64×99=", "15×59="),
    @("33×55=", "38×63="),
    @("60×67=", "47×29="),
    @("55×20=", "60×58="),
    @("85×57=", "96×77="),
    @("63×75=", "65×74="),
    @("32×87=", "79×28="),
    @("13×90=", "39×58="),
    @("18×64=", "61×38="),
    @("57×75=", "32×68="),
    @("12×20=", "32×23="),
    @("42×53=", "18×79="),
    @("38×47=", "13×55="),
    @("65×92=", "24×41="),
    @("93×61=", "86×66="),
    @("12×68=", "66×78="),
    @("89×16=", "65×21="),
    @("22×91=", "36×69="),
    @("55×79=", "60×21="),
    @("16×56=", "80×40="),
    @("61×83=", "15×65="),
    @("95×92=", "17×47="),
    @("76×38=", "27×34=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
